$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text format to column D so numeric-looking strings stay as text
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '24.445.30'
$ws.Range("E2").Value = '  -0.15%  '

# Row 3
$ws.Range("D3").Value = '1.655.54'
$ws.Range("E3").Value = '  -2.57%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -1.17%  '

# Row 5
$ws.Range("D5").Value = '307.77'
$ws.Range("E5").Value = '  -0.61%  '

# Row 6
$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.74%  '

# Row 7
$ws.Range("D7").Value = '0.3627'
$ws.Range("E7").Value = '  -2.99%  '

# Row 8
$ws.Range("D8").Value = '47.36'
$ws.Range("E8").Value = '  -3.72%  '

# Row 9
$ws.Range("D9").Value = '0.3267'
$ws.Range("E9").Value = '  -5.00%  '

# Row 10
$ws.Range("D10").Value = '1.126'
$ws.Range("E10").Value = '  -4.06%  '

# Row 11
$ws.Range("D11").Value = '0.06955'
$ws.Range("E11").Value = '  -6.51%  '

# Row 12
$ws.Range("D12").Value = '0.9987'
$ws.Range("E12").Value = '  -1.29%  '

# Row 13
$ws.Range("D13").Value = '5.943'
$ws.Range("E13").Value = '  -4.44%  '

# Row 14
$ws.Range("D14").Value = '19.34'
$ws.Range("E14").Value = '  -6.58%  '

# Row 15
$ws.Range("D15").Value = '6.624'
$ws.Range("E15").Value = '  -3.75%  '

# Row 16
$ws.Range("D16").Value = '1.652.48'
$ws.Range("E16").Value = '  -2.95%  '

# Row 17
$ws.Range("D17").Value = '0.00001043'
$ws.Range("E17").Value = '  -6.24%  '

# Row 18
$ws.Range("D18").Value = '0.06515'
$ws.Range("E18").Value = '  -3.11%  '

# Row 19
$ws.Range("D19").Value = '0.9992'
$ws.Range("E19").Value = '  -0.81%  '

# Row 20
$ws.Range("D20").Value = '76.45'
$ws.Range("E20").Value = '  -8.17%  '

# Row 21
$ws.Range("D21").Value = '5.925'
$ws.Range("E21").Value = '  -6.37%  '

# Row 22
$ws.Range("D22").Value = '15.73'
$ws.Range("E22").Value = '  -7.65%  '

# Row 23
$ws.Range("D23").Value = '12.66'
$ws.Range("E23").Value = '  -1.59%  '

# Row 24
$ws.Range("D24").Value = '24.389.83'
$ws.Range("E24").Value = '  -0.59%  '

# Row 25
$ws.Range("D25").Value = '2.448'
$ws.Range("E25").Value = '  -0.46%  '

# Row 26
$ws.Range("D26").Value = '2.343'
$ws.Range("E26").Value = '  -14.48%  '

# Row 27
$ws.Range("D27").Value = '146.09'
$ws.Range("E27").Value = '  -2.53%  '

# Row 28
$ws.Range("D28").Value = '18.41'
$ws.Range("E28").Value = '  -8.75%  '

# Row 29
$ws.Range("D29").Value = '1.836.42'
$ws.Range("E29").Value = '  -3.16%  '

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '1.212'
$ws.Range("E30").Value = '  +4.06%  '

# Row 31
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '124.03'
$ws.Range("E31").Value = '  -5.37%  '

# Row 32
$ws.Range("D32").Value = '4.050'
$ws.Range("E32").Value = '  -4.79%  '

# Row 33
$ws.Range("D33").Value = '5.619'
$ws.Range("E33").Value = '  -15.88%  '

# Row 34
$ws.Range("D34").Value = '0.08326'
$ws.Range("E34").Value = '  -4.95%  '

# Row 35
$ws.Range("D35").Value = '1.673'
$ws.Range("E35").Value = '  -6.48%  '

# Row 36
$ws.Range("D36").Value = '12.34'
$ws.Range("E36").Value = '  -8.91%  '

# Row 37
$ws.Range("D37").Value = '5.212'
$ws.Range("E37").Value = '  -5.74%  '

# Row 38
$ws.Range("D38").Value = '0.06054'
$ws.Range("E38").Value = '  -6.58%  '

# Row 39
$ws.Range("D39").Value = '0.02204'
$ws.Range("E39").Value = '  -7.16%  '

# Row 40
$ws.Range("D40").Value = '1.206'
$ws.Range("E40").Value = '  -5.07%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '8.200'
$ws.Range("E41").Value = '  -7.70%  '

# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.2048'
$ws.Range("E42").Value = '  -5.96%  '

# Row 43
$ws.Range("D43").Value = '0.9988'
$ws.Range("E43").Value = '  -0.85%  '

# Row 44
$ws.Range("D44").Value = '0.5850'
$ws.Range("E44").Value = '  -8.02%  '

# Row 45
$ws.Range("D45").Value = '3.734'
$ws.Range("E45").Value = '  -2.23%  '

# Row 46
$ws.Range("D46").Value = '12.58'
$ws.Range("E46").Value = '  -9.14%  '

# Row 47
$ws.Range("D47").Value = '0.5578'
$ws.Range("E47").Value = '  -7.53%  '

# Row 48
$ws.Range("D48").Value = '122.17'
$ws.Range("E48").Value = '  -4.90%  '

# Row 49
$ws.Range("D49").Value = '1.938'
$ws.Range("E49").Value = '  -8.09%  '

# Row 50
$ws.Range("D50").Value = '0.06904'
$ws.Range("E50").Value = '  -4.54%  '

# Row 51
$ws.Range("D51").Value = '74.06'
$ws.Range("E51").Value = '  -6.46%  '

# Reset style on column D back to Normal (removes temporary text-format stamp)
$ws.Range("D2:D51").Style = "Normal"